$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, shifting existing rows (2-4) down to (3-5)
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with slugified identifiers that relate to the
# column headers in row 1, enabling SKOS hierarchical relations between
# columns (fixes issue #13).
$ws.Range("A2").Value = "poblacion"
$ws.Range("B2").Value = "municipio-codigo"
$ws.Range("C2").Value = "tipo-de-estudios-realizados"
$ws.Range("D2").Value = "tipo-de-estudios-realizados-codigo"
$ws.Range("E2").Value = "municipio-nombre"
